# Backup QR Scanner data - append new scan/manual log rows and rename sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Session" to "Pharmacology"
$ws.Name = "Pharmacology"

# New log rows to append (Student ID, Subject, Log Date, Log Time, Type, User)
$rows = @(
    @("244647", "Pharmacology", "14/10/2025", "08:44:52", "Scan",   "marian.samir@med.asu.edu.eg"),
    @("244910", "Pharmacology", "14/10/2025", "08:45:12", "Manual", "marian.samir@med.asu.edu.eg"),
    @("244783", "Pharmacology", "14/10/2025", "08:45:17", "Manual", "marian.samir@med.asu.edu.eg")
)

$startRow = 5
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # Column A holds a numeric-looking Student ID that must be kept as text,
    # exactly like the existing rows (numberStoredAsText). Format the cell
    # as text before assigning so Excel doesn't coerce it into a number.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[0]

    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
